$excel.DisplayAlerts = $false
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Sheet housekeeping: delete the empty "Sheet1", rename "errorMessages" to
#    "errorInfoMessages", and insert a brand-new "credentialsErrorMessages"
#    sheet right after it.
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("Sheet1").Delete()

$errSheet = $wb.Worksheets.Item("errorMessages")
$errSheet.Name = "errorInfoMessages"

$credErrSheet = $wb.Worksheets.Add($null, $errSheet)
$credErrSheet.Name = "credentialsErrorMessages"

# ---------------------------------------------------------------------------
# 2. Add a new row to the "credentials" sheet with the invalid-login test
#    data (invalid / hola_mundo / adios_mundo).
# ---------------------------------------------------------------------------
$credSheet = $wb.Worksheets.Item("credentials")

$credSheet.Range("A6").Value = "invalid"
$credSheet.Range("B6").Value = "hola_mundo"
$credSheet.Range("C6").Value = "adios_mundo"

$bodyRange = $credSheet.Range("A6:C6")
$bodyRange.NumberFormat = "@"
$bodyRange.Borders().LineStyle = 1
$bodyRange.Borders().Weight = 2
$bodyRange.Borders().Color = 0

# Header cells B1:C1 lose their bottom border (the table now continues below).
$headerRange = $credSheet.Range("B1:C1")
$hLeft = $headerRange.Borders().Item(7)   # xlEdgeLeft
$hLeft.LineStyle = 1
$hLeft.Weight = 2
$hLeft.Color = 0

$hTop = $headerRange.Borders().Item(8)    # xlEdgeTop
$hTop.LineStyle = 1
$hTop.Weight = 2
$hTop.Color = 0

$hRight = $headerRange.Borders().Item(10) # xlEdgeRight
$hRight.LineStyle = 1
$hRight.Weight = 2
$hRight.Color = 0

$hBottom = $headerRange.Borders().Item(9) # xlEdgeBottom
$hBottom.LineStyle = -4142                # xlLineStyleNone

# ---------------------------------------------------------------------------
# 3. Populate the new "credentialsErrorMessages" sheet.
# ---------------------------------------------------------------------------
$credErrSheet.Range("A1").Value = "key"
$credErrSheet.Range("B1").Value = "message"
$credErrSheet.Range("A2").Value = "invalid"
$credErrSheet.Range("B2").Value = "Username and password do not match any user in this service."
$credErrSheet.Range("A3").Value = "locked"
$credErrSheet.Range("B3").Value = "Sorry, this user has been locked out."

$credErrSheet.Columns.Item(2).ColumnWidth = 46.83203125

$credErrHeader = $credErrSheet.Range("A1:B1")
$credErrHeader.Font.Bold = $true
$credErrHeader.Interior.Color = 65535
$credErrHeader.HorizontalAlignment = -4108  # xlCenter
$credErrHeader.Borders().LineStyle = 1
$credErrHeader.Borders().Weight = 2

$credErrBody = $credErrSheet.Range("A2:B3")
$credErrBody.Borders().LineStyle = 1
$credErrBody.Borders().Weight = 2

# ---------------------------------------------------------------------------
# 4. Cosmetic window-state tweaks (grid lines / selections / active tab).
# ---------------------------------------------------------------------------
$itemDataSheet = $wb.Worksheets.Item("itemData")
$itemDataSheet.Activate()
$excel.ActiveWindow.DisplayGridlines = $false

$errSheet.Activate()
$excel.ActiveWindow.DisplayGridlines = $false
$errSheet.Range("B11").Select()

$credSheet.Activate()
$excel.ActiveWindow.DisplayGridlines = $false
$credSheet.Range("B6").Select()

$urlsSheet = $wb.Worksheets.Item("urls")
$urlsSheet.Activate()
$excel.ActiveWindow.DisplayGridlines = $false
$urlsSheet.Range("B9").Select()

$credErrSheet.Activate()
$excel.ActiveWindow.DisplayGridlines = $false
$credErrSheet.Range("A3").Select()

Write-Host "done"
